$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated orders: distance codes and one size code were renumbered.
# Apply as global text substitutions across the sheet (mirrors Excel's
# Find & Replace, which rewrites every matching shared string / cell in
# one pass without disturbing row/column layout).
$ws.Cells.Replace("D64", "D69")
$ws.Cells.Replace("D80", "D86")
$ws.Cells.Replace("D51", "D55")
$ws.Cells.Replace("S30", "S31")
